$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.558.15"
$ws.Range('E2').Value = "'  -2.42%  "
$ws.Range('D3').Value = "'1.816.13"
$ws.Range('E3').Value = "'  -2.02%  "
$ws.Range('D4').Value = "'1.007"
$ws.Range('E4').Value = "'  +0.62%  "
$ws.Range('D5').Value = "'308.96"
$ws.Range('E5').Value = "'  -1.58%  "
$ws.Range('E6').Value = "'  +0.64%  "
$ws.Range('D7').Value = "'0.4569"
$ws.Range('E7').Value = "'  -1.39%  "
$ws.Range('D8').Value = "'0.3666"
$ws.Range('E8').Value = "'  -1.24%  "
$ws.Range('D9').Value = "'0.07135"
$ws.Range('E9').Value = "'  -2.07%  "
$ws.Range('D10').Value = "'0.8781"
$ws.Range('E10').Value = "'  -0.96%  "
$ws.Range('D11').Value = "'0.07787"
$ws.Range('E11').Value = "'  -1.04%  "
$ws.Range('D12').Value = "'19.37"
$ws.Range('E12').Value = "'  -3.63%  "
$ws.Range('D13').Value = "'1.825.19"
$ws.Range('E13').Value = "'  -2.45%  "
$ws.Range('D14').Value = "'5.294"
$ws.Range('E14').Value = "'  -1.81%  "
$ws.Range('D15').Value = "'6.376"
$ws.Range('E15').Value = "'  -2.12%  "
$ws.Range('D16').Value = "'86.31"
$ws.Range('E16').Value = "'  -5.20%  "
$ws.Range('D17').Value = "'1.009"
$ws.Range('E17').Value = "'  +0.74%  "
$ws.Range('D18').Value = "'0.000008612"
$ws.Range('E18').Value = "'  -3.44%  "
$ws.Range('E19').Value = "'  +0.65%  "
$ws.Range('D20').Value = "'26.626.37"
$ws.Range('E20').Value = "'  -2.27%  "
$ws.Range('D21').Value = "'14.28"
$ws.Range('E21').Value = "'  -2.72%  "
$ws.Range('D22').Value = "'5.006"
$ws.Range('E22').Value = "'  -1.52%  "
$ws.Range('D23').Value = "'10.47"
$ws.Range('E23').Value = "'  -0.48%  "
$ws.Range('D24').Value = "'1.987"
$ws.Range('E24').Value = "'  +1.92%  "
$ws.Range('D25').Value = "'151.69"
$ws.Range('E25').Value = "'  +0.26%  "
$ws.Range('E26').Value = "'  -2.25%  "
$ws.Range('D27').Value = "'2.076"
$ws.Range('E27').Value = "'  +1.35%  "
$ws.Range('D28').Value = "'113.12"
$ws.Range('E28').Value = "'  -2.48%  "
$ws.Range('D29').Value = "'4.863"
$ws.Range('E29').Value = "'  -3.55%  "
$ws.Range('D30').Value = "'0.08694"
$ws.Range('E30').Value = "'  -1.23%  "
$ws.Range('D31').Value = "'3.063"
$ws.Range('E31').Value = "'  -2.45%  "
$ws.Range('D32').Value = "'4.517"
$ws.Range('E32').Value = "'  -0.06%  "
$ws.Range('D33').Value = "'0.7351"
$ws.Range('E33').Value = "'  -4.47%  "
$ws.Range('B34').Value = "'ARBITRUM"
$ws.Range('C34').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('D34').Value = "'1.119"
$ws.Range('E34').Value = "'  -3.94%  "
$ws.Range('B35').Value = "'RenderToken"
$ws.Range('C35').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('D35').Value = "'2.672"
$ws.Range('E35').Value = "'  -2.05%  "
$ws.Range('E36').Value = "'  -0.04%  "
$ws.Range('D37').Value = "'1.083"
$ws.Range('E37').Value = "'  -2.40%  "
$ws.Range('D38').Value = "'0.01949"
$ws.Range('E38').Value = "'  +0.52%  "
$ws.Range('B39').Value = "'MXToken"
$ws.Range('C39').Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range('D39').Value = "'2.912"
$ws.Range('E39').Value = "'  -1.00%  "
$ws.Range('B40').Value = "'Hedera"
$ws.Range('C40').Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range('D40').Value = "'0.05120"
$ws.Range('E40').Value = "'  -1.90%  "
$ws.Range('D41').Value = "'6.999"
$ws.Range('E41').Value = "'  -0.57%  "
$ws.Range('D42').Value = "'0.5024"
$ws.Range('E42').Value = "'  -1.97%  "
$ws.Range('D43').Value = "'0.1559"
$ws.Range('E43').Value = "'  -4.28%  "
$ws.Range('D44').Value = "'8.187"
$ws.Range('E44').Value = "'  -3.17%  "
$ws.Range('D45').Value = "'1.008"
$ws.Range('E45').Value = "'  +0.78%  "
$ws.Range('D46').Value = "'0.4627"
$ws.Range('E46').Value = "'  -3.60%  "
$ws.Range('D47').Value = "'9.976"
$ws.Range('E47').Value = "'  -3.87%  "
$ws.Range('D48').Value = "'101.05"
$ws.Range('D49').Value = "'1.595"
$ws.Range('E49').Value = "'  -2.94%  "
$ws.Range('D50').Value = "'0.06006"
$ws.Range('E50').Value = "'  -3.16%  "
$ws.Range('D51').Value = "'64.51"
$ws.Range('E51').Value = "'  -1.23%  "
